$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A3: update the clinical-trial "Projet" reference embedded in the big
#         info block, keeping the rest of the text identical. ---
$infoText = $ws.Range("A3").Value2
$infoText = $infoText -replace "Projet: P-2021A22L", "Projet: P94C001FR"
$ws.Range("A3").Value = $infoText

# --- D3: the "Date" observation value switches from a free-text timestamp
#         to a real date/time serial, displayed as yyyy-mm-dd hh:mm:ss.
#         Setting NumberFormat before Value makes the engine store it as a
#         genuine number instead of re-using the cell's old text type. ---
$dateRange = $ws.Range("D3:E3")
$dateRange.NumberFormat = "yyyy-mm-dd hh:mm:ss"
$dateRange.HorizontalAlignment = -4152
$ws.Range("D3").Value = 44326.615185185183

# --- D6 / D7: "Effets secondaires" / "Effets secondaires graves" values ---
$ws.Range("D6").Value = "0.2"
$ws.Range("D7").Value = "0"
